$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert a duplicate of A1 into B1, shifting Noms et
# prenoms / Numero du role / Pseudo de l'eleve one column to the right
# (B1->C1, C1->D1, D1->E1), so the header row reads:
# A1=B1="Numeros d'admission", C1="Noms et prenoms",
# D1="Numero du role", E1="Pseudo de l'eleve"
$oldD1 = $ws.Range("D1").Value2
$oldC1 = $ws.Range("C1").Value2
$oldB1 = $ws.Range("B1").Value2
$oldA1 = $ws.Range("A1").Value2

$ws.Range("E1").Value = $oldD1
$ws.Range("D1").Value = $oldC1
$ws.Range("C1").Value = $oldB1
$ws.Range("B1").Value = $oldA1

# --- Row 4 edits ---
# Student name typo fix: "bobette " -> "zbobette "
$ws.Range("B4").Value = "zbobette "

# Numeric tweaks
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 31

# K4 / L4 flip from "true" to "false" (text values, not booleans) -
# copy from I4, which already holds the text "false", so we don't
# introduce a native boolean type or a stray quote-prefix style.
$ws.Range("I4").Copy($ws.Range("K4"))
$ws.Range("I4").Copy($ws.Range("L4"))

# --- Row 26 edit ---
$ws.Range("G26").Value = 20
